# Apply textual corrections: generalize "CivicActions" possessives and
# tidy up wording, per the commit "New 'Contractor' component
# (generalized 'CivicActions'); use secrender against templates".

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false,
                             $true, 1, $false, $replace, 2) | Out-Null
}

# 1. CA control: "CivicActions systems" -> "CivicActions' systems" (x2)
Replace-Text `
    "This control is not applicable. CivicActions systems do not have system interconnections. The only communication conducted to CivicActions systems is through the Internet." `
    "This control is not applicable. CivicActions' systems do not have system interconnections. The only communication conducted to CivicActions' systems is through the Internet."

# 2. "is real time" -> "is carried out in real-time" (text of the paragraph's
#    first run; the middle run is just a literal space and is left alone)
Replace-Text `
    "Configuration management and log analysis is real time. OpenSCAP security scans are performed and reviewed monthly. See also: RA-5 and SI-4." `
    "Configuration management and log analysis is carried out in real-time. OpenSCAP security scans are performed and reviewed monthly. See also: RA-5 and SI-4."

# 3. "CivicActions Security." -> "CivicActions' Security Office." (text of the
#    paragraph's third/last run)
Replace-Text `
    "Quarterly review of the control assessments supporting the monitoring is conducted by CivicActions Operations in collaboration with CivicActions Security." `
    "Quarterly review of the control assessments supporting the monitoring is conducted by CivicActions Operations in collaboration with CivicActions' Security Office."

# 4. "security related" -> "security-related"
Replace-Text `
    "CivicActions Security reviews the results of the security scans and security assessments with associated JIRA and/or GitLab Issue tickets created to correlate and analyze security related information generated from the monitoring tools becoming POA&M items for tracking." `
    "CivicActions Security reviews the results of the security scans and security assessments with associated JIRA and/or GitLab Issue tickets created to correlate and analyze security-related information generated from the monitoring tools becoming POA&M items for tracking."

# 5. "though JIRA" -> "through JIRA"; "Information included" -> "The information included"
Replace-Text `
    "POA&M items are tracked by CivicActions Security though JIRA tickets with a security categorization assigned. Information included in the POA&M item include the severity, the due date, the weakness source identifier, and the plugin ID that identified the vulnerability." `
    "POA&M items are tracked by CivicActions Security through JIRA tickets with a security categorization assigned. The information included in the POA&M item include the severity, the due date, the weakness source identifier, and the plugin ID that identified the vulnerability."

# 6. "CivicActions Security to be reviewed" -> "CivicActions' Security Office to be reviewed"
Replace-Text `
    "The security status of the system is reported up to the System Owner and Project Manager via CivicActions Security to be reviewed alongside other security issues relating to the system." `
    "The security status of the system is reported up to the System Owner and Project Manager via CivicActions' Security Office to be reviewed alongside other security issues relating to the system."
